$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "dim"
$ws.Range("F4").Value = "URI-Municipio"
